# "Added New Mac-Address and Document Types"
# Appends 5 new reg_center_machine_device rows (regcntr_id 10002 / machine_id
# 10032, device_ids 3000176-3000180) below the existing data, mirroring the
# shape of the rows already present (lang_code "eng", is_active TRUE,
# cr_by "superadmin", cr_dtimes "now()").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns: A regcntr_id | B machine_id | C device_id | D lang_code
#          E is_active  | F cr_by      | G cr_dtimes
$regcntrId = 10002
$machineId = 10032
$deviceIds = @(3000176, 3000177, 3000178, 3000179, 3000180)
$firstNewRow = 157

for ($i = 0; $i -lt $deviceIds.Length; $i++) {
    $row = $firstNewRow + $i

    $ws.Cells.Item($row, 1).Value = $regcntrId
    $ws.Cells.Item($row, 2).Value = $machineId
    $ws.Cells.Item($row, 3).Value = $deviceIds[$i]
    $ws.Cells.Item($row, 4).Value = "eng"
    $ws.Cells.Item($row, 5).Value = $true
    $ws.Cells.Item($row, 6).Value = "superadmin"
    $ws.Cells.Item($row, 7).Value = "now()"
}

$lastNewRow = $firstNewRow + $deviceIds.Length - 1

# Match the author's resulting scroll/selection state as closely as possible.
$excel.ActiveWindow.ScrollRow = 150
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E" + $lastNewRow).Select()
